$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 83
$ws1.Range("F7").Value = 2654
$ws1.Range("F9").Value = 242
$ws1.Range("F10").Value = 101
$ws1.Range("F11").Value = 9743
$ws1.Range("F15").Value = 11664
$ws1.Range("F16").Value = 11949
$ws1.Range("F18").Value = 83

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 83
$ws4.Range("F7").Value = 2654
$ws4.Range("F10").Value = 242
$ws4.Range("F11").Value = 101
$ws4.Range("F12").Value = 9743
$ws4.Range("F16").Value = 11664
$ws4.Range("F17").Value = 11949
$ws4.Range("F19").Value = 83
